$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: price update
$ws.Range("A3").Value = 1.519

# Row 5: price update
$ws.Range("A5").Value = 1.679

# Row 7: price + station data update
$ws.Range("A7").Value = 1.709
$ws.Range("B7").Value = "REPSOL"
$ws.Range("C7").Value = "CL MADRID, 52"

# Row 8: address swap (price/locality unchanged)
$ws.Range("C8").Value = "CALLE COPENHAGUES/N, S/N"

# Row 9: address swap
$ws.Range("C9").Value = "A-6 km 25,5"

# Row 10: address swap
$ws.Range("C10").Value = "CTRA. M-505 km 5,5"

# Row 11: address swap
$ws.Range("C11").Value = "CARRETERA M-505 km 5.5"

# Row 12: price + station data update
$ws.Range("A12").Value = 1.719
$ws.Range("B12").Value = "BP VALDONAIRE"
$ws.Range("C12").Value = "CARRETERA AVD.DE LA INDUSTRIA KM. 15"

# Row 13: price + station data update
$ws.Range("A13").Value = 1.719
$ws.Range("B13").Value = "BP HUMANES - EL MOLINO"
$ws.Range("C13").Value = "AVENIDA DE LAS FLORES, 2"

# Row 14: price + station data update
$ws.Range("A14").Value = 1.719
$ws.Range("B14").Value = "CEPSA"
$ws.Range("C14").Value = "CARRETERA M-405 KM. 5,6"

# Row 15: price + station data update
$ws.Range("A15").Value = 1.719
$ws.Range("B15").Value = "REPSOL HUMANES"
$ws.Range("C15").Value = "AVENIDA LA INDUSTRIA, S/N"

# Row 16: price + station data update (also locality changes)
$ws.Range("A16").Value = 1.719
$ws.Range("B16").Value = "REPSOL"
$ws.Range("C16").Value = "CARRETERA AVENIDA  DE LA INDUSTRIA , 46 KM. 1,1"
$ws.Range("D16").Value = "HUMANES DE MADRID"

# Row 17: price + station data update
$ws.Range("A17").Value = 1.719
$ws.Range("B17").Value = "BP LAS ROZAS"
$ws.Range("C17").Value = "CL LAS CRUCES  S/N"

# Row 18: price + station data update (also locality changes)
$ws.Range("A18").Value = 1.719
$ws.Range("B18").Value = "REPSOL"
$ws.Range("C18").Value = "CR A-6, 20,3"
$ws.Range("D18").Value = "ROZAS DE MADRID (LAS)"

# Row 19: price + station data update
$ws.Range("A19").Value = 1.725
$ws.Range("B19").Value = "CEPSA"
$ws.Range("C19").Value = "CARRETERA M-405 KM. 6"
